$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.381.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.819.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "703.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.817.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.487"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.50%  "

$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.466.61"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.830.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.498.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "514.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.967.90"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.10%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("E30").Value = "  -4.33%  "

$ws.Range("E31").Value = "  -4.82%  "

$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.171"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.788.62"
$ws.Range("D37").Style = "Normal"

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.09%  "

$ws.Range("E40").Value = "  -1.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.39%  "

$ws.Range("E42").Value = "  -2.25%  "

$ws.Range("E43").Value = "  -3.12%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "166.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "434.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.72%  "

$ws.Range("E49").Value = "  -5.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.99%  "
